# Plantilla de Casos de Uso - "Modulo de catalogos terminado" edit
#
# This script updates the "Estado" (Status) and " Esfuerzo (hrs)" (Effort)
# columns for several Use Case rows on the "Casos de Uso" sheet:
#   - Several rows move from "planificado" to "En proceso", with their
#     effort estimate increased.
#   - CU-25 (row 29) moves from the (legacy/empty) "vacio" status to
#     "planificado".
#   - CU-26 (row 30) moves from "vacio" to "En proceso" and gets an
#     effort estimate.
#
# It also updates the sheet's active selection/scroll position to match
# the author's final cursor position (cell E30, scrolled so row 11 is on
# top), since a view change is part of the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Column layout on this sheet:
#   B = Identificador (ID) de CU
#   C = Descripción de Caso de Uso
#   D = Alias
#   E = Estado
#   F =  Esfuerzo (hrs)
#   G = Incremento
#   H = Prioridad
#   I = Comentarios

# Row -> (Estado, Esfuerzo)
$updates = @{
    13 = @("En proceso", 3)
    14 = @("En proceso", 3)
    24 = @("En proceso", 4)
    25 = @("En proceso", 3)
    26 = @("En proceso", 4)
    27 = @("En proceso", 3)
    28 = @("En proceso", 3)
    29 = @("planificado", 0)
    30 = @("En proceso", 2)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("E$row").Value = $vals[0]
    $ws.Range("F$row").Value = $vals[1]
}

# Update the active view/selection to match the saved cursor state.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E30").Select()
